# sampling cohort 6 2018-08-21
# Fill in the Site_name (column C) for the new batch of samples on the
# Data_collect sheet, and leave the view focused on that sheet/range.

$wb = $excel.ActiveWorkbook
$wsData = $wb.Worksheets.Item("Data_collect")

# Column C ("Site_name") values for rows 2-91, alternating in blocks of 15
# between the "OG" and "Ferris" collection sites.
$wsData.Range("C2:C16").Value = "OG"
$wsData.Range("C17:C31").Value = "Ferris"
$wsData.Range("C32:C46").Value = "OG"
$wsData.Range("C47:C61").Value = "Ferris"
$wsData.Range("C62:C76").Value = "OG"
$wsData.Range("C77:C91").Value = "Ferris"

# Switch the active/visible sheet to Data_collect and set the selection
# to reflect where the user left off editing.
$wsData.Activate() | Out-Null
$wsData.Range("C77").Select() | Out-Null
